# "bug fixes and city area" - append two new entries to the work log
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: bug fixes entry
$ws.Range("A43").Value = 45680
$ws.Range("A43").NumberFormat = "d-mmm"
$ws.Range("B43").Value = "butfixes"
$ws.Range("C43").Value = 4

# Row 44: new ideas / city area entry
$ws.Range("A44").Value = 45681
$ws.Range("A44").NumberFormat = "d-mmm"
$ws.Range("B44").Value = "testing out some new ideas"
$ws.Range("C44").Value = 6

# Restore the window to a maximized-like size/position and scroll/select
# the view near the new rows, matching the author's final window state.
$win = $excel.ActiveWindow
$win.Top = -120
$win.Left = -120
$win.Width = 29040
$win.Height = 16440

$win.ScrollRow = 17
$win.ScrollColumn = 1

$ws.Range("C45").Select()
